# SAL-137: Split the single sfGroupsData sheet into three test-data sheets
# ("Edit", "Create", "Delete"), update the group-name fixtures, and add a
# new "Edit" group row with description/font metadata.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. The original "Sheet1" becomes "Create" (keeps its data/sheetId).
# ---------------------------------------------------------------------
$create = $wb.Worksheets.Item(1)
$create.Name = "Create"

# Update existing group-name fixtures (Haseena14/15/16 -> 20/21/22)
$create.Range("A2").Value = "Haseena20"
$create.Range("A3").Value = "Haseena21"
$create.Range("A4").Value = "Haseena22"

# Add a new fifth row re-using the same sample image path as the others
$create.Range("A5").Value = "Haseena23"
$create.Range("B5").Value = "C:\Users\amhas\Pictures\sample images\sample1.jpg"

# ---------------------------------------------------------------------
# 2. Add the "Delete" sheet right after "Create" (consumes sheetId 2 so
#    that the next sheet added before Create picks up sheetId 3, matching
#    the target workbook.xml sheetId layout: Edit=3, Create=1, Delete=4).
# ---------------------------------------------------------------------
$deleteSheet = $wb.Worksheets.Add($null, $wb.Worksheets.Item("Create"))
$deleteSheet.Name = "Delete"

$deleteSheet.Range("A1").Value = "strGroupName"
$deleteSheet.Range("A2").Value = "Haseena22"

$deleteSheet.Range("A1").Interior.Color = 65535
$deleteSheet.Columns.Item(1).ColumnWidth = 11.666666666666666

# ---------------------------------------------------------------------
# 3. Add the "Edit" sheet before "Create" (first tab).
# ---------------------------------------------------------------------
$editSheet = $wb.Worksheets.Add($wb.Worksheets.Item("Create"))
$editSheet.Name = "Edit"

$editSheet.Range("A1").Value = "strGroupName"
$editSheet.Range("B1").Value = "strGroupDescription"
$editSheet.Range("C1").Value = "strFontName"
$editSheet.Range("D1").Value = "strFontSize"

$editSheet.Range("A2").Value = "Haseena21"
$editSheet.Range("B2").Value = "TestLeaf"
$editSheet.Range("C2").Value = "Verdana"
$editSheet.Range("D2").Value = 18

$editSheet.Range("A1:D1").Interior.Color = 65535
$editSheet.Columns.Item(1).ColumnWidth = 11.666666666666666
$editSheet.Columns.Item(2).ColumnWidth = 15.916666666666666

# ---------------------------------------------------------------------
# 4. Selections per sheet (mirrors the saved view state of each tab).
# ---------------------------------------------------------------------
$create.Range("B3:B5").Select()
$editSheet.Range("C3").Select()
$deleteSheet.Range("A3").Select()

# "Delete" ends up the active tab (last-selected sheet).
